# Applies the Travis County 2016 bg SVI refactor-variance refresh:
#   - Sheet "Significant Components": re-order the factor-member lists in
#     column C for the F2 factors (rows 2, 4, 5, 6) to reflect the new
#     significance ordering.
#   - Sheet "Loading Factors": re-order / refresh the F2 loading rows
#     (labels in column A, loadings in B:F) to match the refit.
#   - Sheet "All Refactor Variances": refresh the F1 (I:M) and F2 (N:R)
#     SS Loadings / Proportion / Cumulative / Ratio Variance rows.
#   - Sheet "Final Variances": refresh the F2 (B:F) SS Loadings /
#     Proportion / Cumulative / Ratio Variance rows (mirrors the F2
#     columns of "All Refactor Variances").
#   - Sheet "Included and Excluded": refresh the ordering of the include
#     list in B2 to match the new factor-member ordering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Significant Components
# ---------------------------------------------------------------------
$wsSig = $wb.Worksheets.Item("Significant Components")

$wsSig.Range("C2").Value = "['QSERV' 'QHISPC' 'QEDLESHI' 'QNOHLTH' 'QESL' 'PPUNIT' 'QEXTRCT' 'QFHH'`n 'PERCAP']"
$wsSig.Range("C4").Value = "['QAGEDEP' 'QSSBEN' 'MEDAGE']"
$wsSig.Range("C5").Value = "['QNOAUTO' 'QPOVTY' 'QRENTER']"
$wsSig.Range("C6").Value = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"

# ---------------------------------------------------------------------
# 2) Loading Factors (rows 2-20: column A label, columns B-F = F2:1..F2:5)
# ---------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Loading Factors")

$loadRows = @(
    @(2, "QSERV", 0.527773798401461, 0.3860111216432911, -0.1868758135116943, 0.354321523853212, -0.03855995772938289),
    @(3, "QHISPC", 0.8335417236436697, 0.3477459721311119, -0.1159861013261645, 0.1412652469771052, -0.09805697665199882),
    @(4, "QEDLESHI", 0.8755563575228159, 0.2475053926510344, 0.01400336487767863, 0.2099414451999919, -0.1400362633202189),
    @(5, "QNOHLTH", 0.6484773974648308, 0.4451428969732948, -0.09382909943578427, 0.3175828013834935, -0.1542397774428718),
    @(6, "QESL", 0.7907959742290726, 0.1697374991404037, -0.01961143171220127, 0.2266196868953843, -0.2859865796627279),
    @(7, "PPUNIT", 0.7835523160839951, -0.003789135076196481, -0.131320794436579, -0.377247336061985, 0.1043526311122389),
    @(8, "QEXTRCT", 0.7523372406717769, 0.1476097716879644, -0.01125724233229494, 0.09236455865600428, -0.2688479104752657),
    @(9, "QFHH", 0.5447678211090614, 0.3097178102215681, -0.08643488404640712, 0.07397525800617698, 0.2884044969508311),
    @(10, "MDHSEVAL", 0.3668007681726036, 0.8225015767674182, -0.02259948959778094, -0.03119506946875829, -0.02428594007996773),
    @(11, "QRICH", 0.1818961232561695, 0.8639279978086627, -0.1646723816573004, 0.3092163333297806, -0.01566263042680391),
    @(12, "PERCAP", 0.4803068718055284, 0.6956540116243626, -0.2531491223390959, 0.2524414130048208, 0.06893419161578869),
    @(13, "QNOAUTO", 0.1274383136779815, 0.08567042951629011, -0.06014622973025791, 0.6906315125289809, 0.04382609677170136),
    @(14, "QPOVTY", 0.3961340608327603, 0.1517075738531815, -0.3133146727667287, 0.5535940069429346, 0.09229548562342259),
    @(15, "QRENTER", -0.01907546957071893, 0.2341572257899979, -0.4268194506734547, 0.7484008082129754, -0.1155640098205284),
    @(16, "QAGEDEP", -0.03648373641856956, -0.1156683492432435, 0.6860459839132048, -0.09660788802146582, 0.5943324587269403),
    @(17, "QSSBEN", 0.01908852566369217, -0.0371254773654064, 0.7831494100350442, -0.1342903110163737, 0.09955596294146797),
    @(18, "MEDAGE", -0.3172320232961239, -0.2425275147491238, 0.7789709694832544, -0.2970161715264348, -0.05364863570250369),
    @(19, "QFEMLBR", -0.2367556356186458, 0.08379889668760601, -0.03449523820275541, 0.04713344811063635, 0.8071208598926645),
    @(20, "QFEMALE", -0.05366633319430788, -0.0638708668468258, 0.18446360544655, -0.006743985974391692, 0.8303066162144787)
)

foreach ($r in $loadRows) {
    $rowNum = $r[0]
    $wsLoad.Cells.Item($rowNum, 1).Value = $r[1]
    $wsLoad.Cells.Item($rowNum, 2).Value = $r[2]
    $wsLoad.Cells.Item($rowNum, 3).Value = $r[3]
    $wsLoad.Cells.Item($rowNum, 4).Value = $r[4]
    $wsLoad.Cells.Item($rowNum, 5).Value = $r[5]
    $wsLoad.Cells.Item($rowNum, 6).Value = $r[6]
}

# ---------------------------------------------------------------------
# 3) All Refactor Variances (rows 2-5: columns I-M = F1:1..F1:5,
#    columns N-R = F2:1..F2:5). Columns A-H (F0) are unchanged.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Refactor Variances")

$allRows = @(
    @(2, 4.722615385164985, 3.480168176063847, 2.298524086378757, 2.143689070833887, 2.018073895317462, 4.995669848601255, 2.752559864241837, 2.18426812341469, 2.119818585562304, 2.038891321407597),
    @(3, 0.2248864469126184, 0.1657222940982784, 0.109453527922798, 0.1020804319444708, 0.09609875691987912, 0.262929992031645, 0.1448715718022019, 0.1149614801797205, 0.1115693992401213, 0.1073100695477683),
    @(4, 0.2248864469126184, 0.3906087410108968, 0.5000622689336948, 0.6021427008781656, 0.6982414577980447, 0.262929992031645, 0.4078015638338469, 0.5227630440135674, 0.6343324432536886, 0.7416425128014569),
    @(5, 0.322075471745568, 0.2373423867166177, 0.156755985627046, 0.1461964637081117, 0.1376296922026566, 0.3545238945896744, 0.1953388179636152, 0.1550092911279704, 0.1504355498967859, 0.1446924464219542)
)

foreach ($r in $allRows) {
    $rowNum = $r[0]
    # Columns I..M = F1, indices 9..13
    $wsAll.Cells.Item($rowNum, 9).Value  = $r[1]
    $wsAll.Cells.Item($rowNum, 10).Value = $r[2]
    $wsAll.Cells.Item($rowNum, 11).Value = $r[3]
    $wsAll.Cells.Item($rowNum, 12).Value = $r[4]
    $wsAll.Cells.Item($rowNum, 13).Value = $r[5]
    # Columns N..R = F2, indices 14..18
    $wsAll.Cells.Item($rowNum, 14).Value = $r[6]
    $wsAll.Cells.Item($rowNum, 15).Value = $r[7]
    $wsAll.Cells.Item($rowNum, 16).Value = $r[8]
    $wsAll.Cells.Item($rowNum, 17).Value = $r[9]
    $wsAll.Cells.Item($rowNum, 18).Value = $r[10]
}

# ---------------------------------------------------------------------
# 4) Final Variances (rows 2-5: columns B-F = F2:1..F2:5)
# ---------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Final Variances")

$finalRows = @(
    @(2, 4.995669848601255, 2.752559864241837, 2.18426812341469, 2.119818585562304, 2.038891321407597),
    @(3, 0.262929992031645, 0.1448715718022019, 0.1149614801797205, 0.1115693992401213, 0.1073100695477683),
    @(4, 0.262929992031645, 0.4078015638338469, 0.5227630440135674, 0.6343324432536886, 0.7416425128014569),
    @(5, 0.3545238945896744, 0.1953388179636152, 0.1550092911279704, 0.1504355498967859, 0.1446924464219542)
)

foreach ($r in $finalRows) {
    $rowNum = $r[0]
    $wsFinal.Cells.Item($rowNum, 2).Value = $r[1]
    $wsFinal.Cells.Item($rowNum, 3).Value = $r[2]
    $wsFinal.Cells.Item($rowNum, 4).Value = $r[3]
    $wsFinal.Cells.Item($rowNum, 5).Value = $r[4]
    $wsFinal.Cells.Item($rowNum, 6).Value = $r[5]
}

# ---------------------------------------------------------------------
# 5) Included and Excluded: refresh the include-list ordering in B2
# ---------------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Included and Excluded")

$wsInc.Range("B2").Value = "[['QSERV', 'QHISPC', 'QEDLESHI', 'QNOHLTH', 'QESL', 'PPUNIT', 'QEXTRCT', 'QFHH', 'PERCAP', 'MDHSEVAL', 'QRICH', 'QAGEDEP', 'QSSBEN', 'MEDAGE', 'QNOAUTO', 'QPOVTY', 'QRENTER', 'QFEMLBR', 'QFEMALE']]"
